$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# "Enterprises density (per 1000 people)" row: update Micro/SMEs/MSMEs values
# (values are stored as text in the sheet, so force text format while writing
# them, then restore the original "Normal" cell style)
$rng = $ws.Range("B13:D13")
$rng.NumberFormat = "@"
$ws.Range("B13").Value = "6.89"
$ws.Range("C13").Value = "1.94"
$ws.Range("D13").Value = "8.83"
$rng.Style = "Normal"
